$p = $ppt.ActivePresentation

# --- 1. Table style swap on slide 5 (graphicFrame with the B1 financial documents table) ---
$s = $p.Slides.Item(5)
$tableShape = $s.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{9BAACAFA-6E49-44FE-9998-D754E6A63700}")

# --- 2. Theme re-colour: swap the "Integral" / Red Violet colour scheme for the
#        "Office Theme" / Office colour scheme (the deck's theme content moved from
#        the Red Violet palette back to the stock Office palette). ---
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = RGB(0x00,0x00,0x00)   # dk1
$colors.Item(2).RGB  = RGB(0xFF,0xFF,0xFF)   # lt1
$colors.Item(3).RGB  = RGB(0x44,0x54,0x6A)   # dk2
$colors.Item(4).RGB  = RGB(0xE7,0xE6,0xE6)   # lt2
$colors.Item(5).RGB  = RGB(0x5B,0x9B,0xD5)   # accent1
$colors.Item(6).RGB  = RGB(0xED,0x7D,0x31)   # accent2
$colors.Item(7).RGB  = RGB(0xA5,0xA5,0xA5)   # accent3
$colors.Item(8).RGB  = RGB(0xFF,0xC0,0x00)   # accent4
$colors.Item(9).RGB  = RGB(0x44,0x72,0xC4)   # accent5
$colors.Item(10).RGB = RGB(0x70,0xAD,0x47)   # accent6
$colors.Item(11).RGB = RGB(0x05,0x63,0xC1)   # hlink
$colors.Item(12).RGB = RGB(0x95,0x4F,0x72)   # folHlink
